# "Preparing my first game." — add a "template" column to the character
# config sheet, between "name" and "moveSpeed".
#
# Before: id | name | moveSpeed | scale | prefab
# After:  id | name | template | moveSpeed | scale | prefab

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting the old moveSpeed/scale/prefab columns
# (C,D,E) right to (D,E,F).
$ws.Columns("C").Insert()

# New header (row1), type (row2), Chinese label (row3) for the template column.
$ws.Range("C1").Value = "template"
$ws.Range("C2").Value = "string"
$ws.Range("C3").Value = "模板"

# Data rows: every character uses the "Character" template for now.
$ws.Range("C4").Value = "Character"
$ws.Range("C5").Value = "Character"

# Match column B's width so the new column looks consistent with its neighbor.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Leave the selection where the last edit was made.
[void]$ws.Range("C5").Select()
